# Append new match row (row 85) to Sheet1, mirroring the format of the
# preceding data row (row 84) and filling in the new match's values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 85

# Clone the formatting (styles) of the last existing data row onto the new
# row before writing values, so the new row reuses the same cell styles
# (bold/bordered index column, date-formatted match-date column) instead of
# creating new style entries.
$ws.Range("A84:V84").Copy()
$ws.Range("A85:V85").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 1).Value = 84
$ws.Cells.Item($row, 2).Value = "thailand"
$ws.Cells.Item($row, 3).Value = "thai-league-1"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45257.52083333334
$ws.Cells.Item($row, 6).Value = "Trat FC"
$ws.Cells.Item($row, 7).Value = 2
$ws.Cells.Item($row, 8).Value = "Khonkaen Utd."
$ws.Cells.Item($row, 9).Value = 2
$ws.Cells.Item($row, 10).Value = 1.89
$ws.Cells.Item($row, 11).Value = "20/11/2023 00:12"
$ws.Cells.Item($row, 12).Value = 1.8
$ws.Cells.Item($row, 13).Value = "27/11/2023 12:29"
$ws.Cells.Item($row, 14).Value = 3.86
$ws.Cells.Item($row, 15).Value = "20/11/2023 00:12"
$ws.Cells.Item($row, 16).Value = 3.93
$ws.Cells.Item($row, 17).Value = "27/11/2023 12:29"
$ws.Cells.Item($row, 18).Value = 3.47
$ws.Cells.Item($row, 19).Value = "20/11/2023 00:12"
$ws.Cells.Item($row, 20).Value = 4.23
$ws.Cells.Item($row, 21).Value = "27/11/2023 12:29"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/trat-fc-khonkaen-united/ETNPm7aO/"
